$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")

$ws.Range("A197").Value = "GSE146931"
$ws.Range("B197").Value = "PRJNA608903"
$ws.Range("C197").Value = "SRR11192628"
$ws.Range("D197").Value = "SRR11192627"
$ws.Range("E197").Value = "single"
$ws.Range("F197").Value = "WIP1"
$ws.Range("H197").Value = "✓"

$ws.Range("A198").Value = "GSE146931"
$ws.Range("B198").Value = "PRJNA608903"
$ws.Range("C198").Value = "SRR11192624"
$ws.Range("D198").Value = "SRR11192623"
$ws.Range("E198").Value = "single"
$ws.Range("F198").Value = "WIP2"
$ws.Range("H198").Value = "✓"

$ws.Range("A199").Value = "GSE146931"
$ws.Range("B199").Value = "PRJNA608903"
$ws.Range("C199").Value = "SRR11192636"
$ws.Range("D199").Value = "SRR11192635"
$ws.Range("E199").Value = "single"
$ws.Range("F199").Value = "WIP2"
$ws.Range("H199").Value = "✓"

$ws.Range("H67").Value = "✓"

$ws.Range("C199").Select()
